$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared string text in A5 from "Changed lightmap 1.21.6+ (Vanilla)"
# to "Changed lightmap 1.21+ (Vanilla)"
$ws.Range("A5").Value = "Changed lightmap 1.21+ (Vanilla)"

# Update the build number in C2 from 7830 to 7840
$ws.Range("C2").Value = 7840

# Move the active selection from C2 to A6
$ws.Range("A6").Select()
